$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.4778628081766
$ws.Range("K2").Value = 13.9286135866664
$ws.Range("M2").Value = 17.3020839614069

$ws.Range("B3").Value = 18.6864511804753
$ws.Range("K3").Value = 16.9554638876458
$ws.Range("M3").Value = 20.8587820982446

$ws.Range("B4").Value = 19.5452760393505
$ws.Range("K4").Value = 25.3385372031742
$ws.Range("M4").Value = 19.8837915557999

$ws.Range("B5").Value = 18.4727295355144
$ws.Range("K5").Value = 18.8097519849247
$ws.Range("M5").Value = 18.7575959940119

$ws.Range("B6").Value = 25.8176804364833
$ws.Range("K6").Value = 24.9676333375889
$ws.Range("M6").Value = 23.1977463905367
